$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.634.33'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.644.61'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.26'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '1.873.40'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('E13').Value = '  +3.28%  '
$ws.Range('D14').Value = '1.638.21'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.532'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').Value = '26.665.36'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '0.0₃0751'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.38'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('E22').Value = '  +1.88%  '
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +10.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.50'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.12'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.88'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0518'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.39'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.06'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').Value = '1.278.22'
$ws.Range('E34').Value = '  +5.45%  '
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('E36').Value = '  +6.25%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.827'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.527'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.76%  '
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('E42').Value = '  -1.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.46'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('D44').Value = '1.784.83'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.12'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.86'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.40%  '
$ws.Range('E47').Value = '  +3.64%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.77'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0978'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.64%  '
